# Update the cryptocurrency Price (D) and Volume(1h) (E) columns
# with freshly scraped values, keeping the cells as plain text so
# Excel does not reinterpret values like '591.17' or '0.460' as numbers/dates.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '63.808.35'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -0.39%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.135.84'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +0.18%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.13%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '591.17'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +0.36%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '145.37'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -1.67%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +0.04%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.128.60'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.08%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.85%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -1.53%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '5.86'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +1.49%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.460'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.83%  '
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -3.63%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.06'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -1.31%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '3.655.31'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.17%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -1.40%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '7.34'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.34%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '3.137.97'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.27%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '63.688.26'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.33%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '469.09'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.35%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.33'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -0.34%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.731'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -0.31%  '
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.67%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.33'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +6.33%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.96'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -2.39%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  -1.18%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.02%  '
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +9.41%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.44'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  +7.95%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.70'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -0.63%  '
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -0.27%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +0.23%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '27.69'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +1.87%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +0.13%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0₃0840'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -5.98%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.66%  '
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -4.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.13'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  +0.27%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.19'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -6.99%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '51.52'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.09%  '
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +6.95%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '452.24'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -0.78%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.292'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +4.77%  '
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -0.82%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.911.34'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.49%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '39.22'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +9.46%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.107'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -3.86%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '131.19'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +4.22%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.24'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +2.19%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.110'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -1.34%  '
